# ---------------------------------------------------------------------------
# Individual test cases will accept multiple groups separated by comma
# values in test runner file.
#
# This script:
#  1. Splits the single "TestRunner" sheet into two sheets:
#       - "Config"     (new, first tab)  -- Properties/Values settings table
#       - "Test Cases" (was "TestRunner", second tab) -- existing test-case grid
#  2. Renames headers Group -> Groups, Execution -> Execute on Test Cases
#  3. Updates the Groups column (and the D2 filter sample) to comma-separated
#     multi-group values
#  4. Re-creates the view/selection state captured in the target workbook
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Config" sheet -------------------------------------
# Worksheets.Add() inserts the new sheet immediately before the (currently)
# active sheet, so it lands in tab position 1 and pushes the original sheet
# to position 2 -- exactly the order/sheetId pattern we need.
$wb.Worksheets.Add() | Out-Null

$wsConfig = $wb.Worksheets.Item(1)
$wsConfig.Name = "Config"

$wsCases = $wb.Worksheets.Item(2)
$wsCases.Name = "Test Cases"

# --- 2. Populate the Config sheet ------------------------------------------
$wsConfig.Range("A1").Value = "Properties"
$wsConfig.Range("B1").Value = "Values"
$wsConfig.Range("A2").Value = "Browser"
$wsConfig.Range("B2").Value = "Chrome"
$wsConfig.Range("A3").Value = "RunInParallel"
$wsConfig.Range("B3").Value = "Yes"
$wsConfig.Range("A4").Value = "NumberOfWindows"
$wsConfig.Range("B4").Value = 4
$wsConfig.Range("A5").Value = "BaseURL"
$wsConfig.Range("B5").Value = "https://vam-bd-agentuw-qa-wapp.azurewebsites.net"

# Column widths
$wsConfig.Columns.Item(1).ColumnWidth = 20
$wsConfig.Columns.Item(2).ColumnWidth = 53

# Borders: row 1 (header) gets a double bottom rule; column A gets a thin
# right rule; column B gets a thin left rule, so the two columns read like a
# bordered two-column table.
$wsConfig.Range("A1:B1").Borders.Item(9).LineStyle = -4119   # xlEdgeBottom / xlDouble
$wsConfig.Range("A1:B1").Borders.Item(9).Color = 0
$wsConfig.Range("A1:A5").Borders.Item(10).LineStyle = 1      # xlEdgeRight / xlThin
$wsConfig.Range("A1:A5").Borders.Item(10).Color = 0
$wsConfig.Range("B1:B5").Borders.Item(7).LineStyle = 1       # xlEdgeLeft / xlThin
$wsConfig.Range("B1:B5").Borders.Item(7).Color = 0

# Alignment: column B (values) and B1 header are left aligned
$wsConfig.Range("B1:B5").HorizontalAlignment = -4131          # xlLeft

$wsConfig.Range("B6").Select()

# --- 3. Update the Test Cases sheet headers ---------------------------------
$wsCases.Range("C1").Value = "Groups"
$wsCases.Range("D1").Value = "Execute"

# --- 4. Update the Groups column to comma-separated multi-group values -----
$wsCases.Range("C2").Value = "Smoke,Smoke"
$wsCases.Range("D2").Value = "Groups=Smoke,Regression"
$wsCases.Range("C3").Value = "Smoke"
$wsCases.Range("C4").Value = "Regression,Sanity"
$wsCases.Range("C5").Value = "Regression"
$wsCases.Range("C6").Value = "Sanity,Regression"
$wsCases.Range("C7").Value = "Sanity,Smoke"
$wsCases.Range("C8").Value = "Sanity"
$wsCases.Range("C9").Value = "Sanity"
$wsCases.Range("C10").Value = "Sanity"
$wsCases.Range("C11").Value = "Sanity"
$wsCases.Range("C12").Value = "Sanity"

$wsCases.Range("C6").Select()

# --- 5. Make "Test Cases" the active tab again (new sheet steals focus) ----
$wb.Worksheets.Item(2).Activate()
